$wb = $excel.ActiveWorkbook

$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")
$wsConfig = $wb.Worksheets.Item("Config")

# --- helper: write a literal text value into a cell without Excel's
#     auto-date/number detection mangling it (e.g. "8/26/2018" -> serial).
#     We do this by writing a formula that evaluates to the literal string,
#     then converting that single cell to a static value via copy/paste-special.
function Set-LiteralText($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# 1. The "when zooming, if a scroll bar is all the way to min or max, keep it
#    there" task (Id 34) is now done: remove it from the Active/Todo sheet...
$wsActive.Rows.Item(2).Delete()

# ...and insert it as a new row at the top of the Inactive/Done sheet, with
# its completion ("Done") date recorded. Clear the formatting Excel copies
# down from the bold header row so the new row matches the plain data rows.
$wsInactive.Rows.Item(2).Insert()
$wsInactive.Rows.Item(2).ClearFormats()
$wsInactive.Cells.Item(2, 1).Value = 34
Set-LiteralText $wsInactive.Cells.Item(2, 2) "when zooming, if a scroll bar is all the way to min or max, keep it there"
$wsInactive.Cells.Item(2, 3).Value = "Done"
$wsInactive.Cells.Item(2, 4).Value = "Task"
Set-LiteralText $wsInactive.Cells.Item(2, 5) "8/11/2018"
Set-LiteralText $wsInactive.Cells.Item(2, 6) "8/26/2018"

# 2. Add the new follow-up task to the Active/Todo sheet (Id 95), inserted
#    right after the current last row (row 11, after the deletion above).
$wsActive.Rows.Item(11).Insert()
$wsActive.Rows.Item(11).ClearFormats()
$wsActive.Cells.Item(11, 1).Value = 95
Set-LiteralText $wsActive.Cells.Item(11, 2) "continue refactoring business logic out of RequestColorWorker"
$wsActive.Cells.Item(11, 3).Value = "Todo"
$wsActive.Cells.Item(11, 4).Value = "Task"
Set-LiteralText $wsActive.Cells.Item(11, 5) "8/26/2018"

# 3. Bump the "Max Id" tracker on the Config sheet to match the newly used Id.
$wsConfig.Cells.Item(2, 6).Value = 95
